$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header label in A1 from "Indicator_code" to "Indicator code"
$ws.Range("A1").Value = "Indicator code"

# Remove the threaded comment that was left on E1 ("Add rows for each year.")
$ws.Range("E1").Comment.Delete()

# Reflect the active cell ending on E1 (where the now-removed comment was)
$null = $ws.Range("E1").Select()
